# Insert a new data row at row 251 (pushing the existing rows 251:361 down
# to 252:362) and populate the new row with the latest weekly price entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(251).Insert()

$ws.Range('A251').Value() = 7
$ws.Range('B251').Value() = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C251').Value() = 'Ñuble'
$ws.Range('D251').Value() = 44917
$ws.Range('E251').Value() = 16
$ws.Range('F251').Value() = 100112008
$ws.Range('G251').Value() = 'Coliflor'
$ws.Range('H251').Value() = 'Sin especificar'
$ws.Range('I251').Value() = 'Segunda'
$ws.Range('J251').Value() = 300
$ws.Range('K251').Value() = 700
$ws.Range('L251').Value() = 750
$ws.Range('M251').Value() = 725
$ws.Range('N251').Value() = '$/unidad'
$ws.Range('O251').Value() = 'Región del Maule'
$ws.Range('P251').Value() = 725
$ws.Range('Q251').Value() = 1
$ws.Range('R251').Value() = 'Hortaliza'

# Match the date cell formatting used throughout column D.
$ws.Range('D251').NumberFormat() = $ws.Range('D252').NumberFormat()
